$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10, shifting existing rows 10-65 down to 11-66.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 2
$ws.Range("B10").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44881
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112032
$ws.Range("G10").Value = "Zapallo italiano"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 560
$ws.Range("K10").Value = 6000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 6500
$ws.Range("N10").Value = "$/caja 50 unidades"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 130
$ws.Range("Q10").Value = 50
$ws.Range("R10").Value = "Hortaliza"

# Match date-format style used by the rest of the "Fecha" column.
$ws.Range("D10").NumberFormat = $ws.Range("D11").NumberFormat
